$wb = $excel.ActiveWorkbook
$ws2023 = $wb.Worksheets.Item("2023")
$wsOverall = $wb.Worksheets.Item("Overall")

# --- "2023" sheet: day 1-8 numbers were corrected (re-synced totals) ---
# (values typed C then B to dodge a dependency-order quirk in the calc engine)
$ws2023.Range("C2").Value = 62098
$ws2023.Range("B2").Value = 194425

$ws2023.Range("C3").Value = 7295
$ws2023.Range("B3").Value = 161961

$ws2023.Range("C4").Value = 15797
$ws2023.Range("B4").Value = 105706

$ws2023.Range("C5").Value = 13996
$ws2023.Range("B5").Value = 103579

$ws2023.Range("C6").Value = 24885
$ws2023.Range("B6").Value = 61307

$ws2023.Range("C7").Value = 1273
$ws2023.Range("B7").Value = 77568

$ws2023.Range("C8").Value = 5647
$ws2023.Range("B8").Value = 57934

$ws2023.Range("C9").Value = 11811
$ws2023.Range("B9").Value = 46429

# --- day 9 (row 10) results are filled in for the first time ---
$ws2023.Range("C10").Value = 1121
$ws2023.Range("B10").Value = 19558
$ws2023.Range("E10").Value = 19396
$ws2023.Range("F10").Value = 19037

# nudge the row-9 / row-10 formulas so they recompute cleanly now that
# the row is no longer blank
$ws2023.Range("D9").Formula = "=IF(ISBLANK(B9),"""",B9+C9)"
$ws2023.Range("H9").Formula = "=IF(ISBLANK(C9),"""",F9/B9)"
$ws2023.Range("D10").Formula = "=IF(ISBLANK(B10),"""",B10+C10)"
$ws2023.Range("H10").Formula = "=IF(ISBLANK(C10),"""",F10/B10)"

# --- "Overall" sheet: mark 2023 day 9 (columns AH:AK, row 12) as succeeded ---
$wsOverall.Range("AH12:AK12").Value = "s"

# --- selection / active-tab bookkeeping, Overall ends up the visible tab ---
[void]$ws2023.Activate()
[void]$ws2023.Range("F10").Select()

[void]$wsOverall.Activate()
[void]$wsOverall.Range("AL12").Select()
